# Update the currency report from "Dolar/USD" (01-02/11/2021) to
# "Coroa Norueguesa/NOK" (11/11/2021), with the new simplified/no-change
# result behaviour described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("C1").Value = "Símbolo"
$ws.Range("D1").Value = "Valor ontem (11/11/2021)"
$ws.Range("E1").Value = "Valor hoje (11/11/2021)"
$ws.Range("F1").Value = "Resultado"

# Data row
$ws.Range("A2").Value = "Coroa Norueguesa"
$ws.Range("B2").Value = "NOK"
$ws.Range("C2").Value = "kr"
$ws.Range("D2").Value = 0.62
$ws.Range("E2").Value = 0.62
$ws.Range("F2").Value = "Nenhuma mudança"

# Report time
$ws.Range("F5").Value = "21:28"
